$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{
    2 = @(1.02, 1.05039458813473, 1.048909380469829, 1.06381790568024, 1.070806667848352, 1.043217386086064, 1.055427972682501, 1.051667911557302, 1.066535518156538, 1.073505522935751, 1.02221899464468)
    3 = @(1.02, 1.051519514112478, 1.04974064279374, 1.064966442646165, 1.072097132856637, 1.043517990724983, 1.056201706364152, 1.052311264445393, 1.067498303623833, 1.074611248611253, 1.022483182104633)
    4 = @(1.02, 1.05224732215546, 1.050278381968585, 1.065709924901993, 1.072932737735765, 1.043711254879945, 1.056701704054577, 1.052726766399158, 1.068121007374703, 1.075326756499781, 1.022653745136543)
    5 = @(1.02, 1.052553271188552, 1.050504413341985, 1.066022557983545, 1.073284167906779, 1.043792204732874, 1.056911745675429, 1.052901254500037, 1.068382725063629, 1.075627564025122, 1.022725357985628)
    6 = @(1.02, 1.052604640150712, 1.050542363002711, 1.066075054766882, 1.073343182933046, 1.043805779076204, 1.05694700337288, 1.05293054075028, 1.068426664709773, 1.075678071429904, 1.022737376713792)
    7 = @(1.02, 1.052251410344378, 1.050281402345035, 1.065714102032086, 1.072937433003937, 1.043712337708124, 1.056704511258831, 1.052729098659656, 1.068124504721478, 1.075330775871977, 1.022654702391766)
    8 = @(1.02, 1.050774782322478, 1.049190339156425, 1.064205997087213, 1.071242665554469, 1.043319235164591, 1.055689596528267, 1.051885499512756, 1.066860955019953, 1.073879202587505, 1.02230835756446)
    9 = @(1.02, 1.048172003276403, 1.047266645538801, 1.061550799110312, 1.068260717138608, 1.042616978712249, 1.053896109992956, 1.05039290285037, 1.06463221938168, 1.071321527139417, 1.021695110378628)
    10 = @(1.02, 1.04643622639997, 1.045983433676909, 1.05978214182167, 1.066275665746778, 1.042142368156681, 1.052697002555117, 1.049393737773433, 1.063144867429654, 1.069616469567295, 1.021284294812087)
    11 = @(1.02, 1.045684458924169, 1.045427608073879, 1.059016630314627, 1.065416789922818, 1.041935326014293, 1.052176948800705, 1.048960110203888, 1.062500453407386, 1.068878160077085, 1.021105935063)
    12 = @(1.02, 1.045405192730569, 1.045221121338485, 1.058732333416308, 1.065097863399774, 1.041858190748276, 1.051983651786805, 1.048798893450753, 1.062261031025098, 1.068603916434536, 1.021039612946004)
    13 = @(1.02, 1.045465097491176, 1.045265414733152, 1.058793313906566, 1.065166269804173, 1.041874746964981, 1.052025120363038, 1.048833481705375, 1.062312390595845, 1.06866274276768, 1.021053842497331)
    14 = @(1.02, 1.045661375238548, 1.045410540396475, 1.05899312926902, 1.065390425363171, 1.041928954695155, 1.05216097338234, 1.048946787002561, 1.062480663874743, 1.068855491073543, 1.021100454312907)
    15 = @(1.02, 1.045782304907997, 1.045499953339014, 1.059116248453277, 1.065528547964502, 1.041962323300638, 1.0522446602037, 1.049016578485511, 1.06258433483419, 1.068974249277322, 1.02112916392934)
    16 = @(1.02, 1.046486115130826, 1.04602031803426, 1.059832953120727, 1.066332680367577, 1.042156076542315, 1.052731499201692, 1.04942249545946, 1.063187626948053, 1.069665468444058, 1.021296121962515)
    17 = @(1.02, 1.046927551520984, 1.046346678906788, 1.060282610174134, 1.066837268136164, 1.042277202264941, 1.053036657009517, 1.049676853017677, 1.06356595341571, 1.070099048769863, 1.021400723332454)
    18 = @(1.02, 1.047185018117824, 1.046537021721288, 1.060544919330234, 1.067131649871905, 1.042347704994518, 1.053214570145959, 1.049825120718676, 1.063786588175138, 1.070351948128884, 1.021461689828827)
    19 = @(1.02, 1.047272804944101, 1.046601920666374, 1.060634365405891, 1.067232037380168, 1.042371719534718, 1.053275220341761, 1.049875660117244, 1.06386181274052, 1.070438180196323, 1.021482470079915)
    20 = @(1.02, 1.046880191217282, 1.04631166533401, 1.06023436294684, 1.066783124012233, 1.042264221916914, 1.053003924771319, 1.049649572664193, 1.063525366366783, 1.070052529787372, 1.021389505324861)
    21 = @(1.02, 1.045603577071934, 1.045367805291474, 1.058934287305515, 1.065324414440352, 1.041912998237963, 1.052120971536306, 1.048913425511204, 1.062431113207089, 1.068798731585098, 1.02108673026635)
    22 = @(1.02, 1.044800766272814, 1.044774198771313, 1.058117156175982, 1.064407833176695, 1.041690835649962, 1.051565094988554, 1.048449723055849, 1.061742774832816, 1.068010403025235, 1.020895950592299)
    23 = @(1.02, 1.045226366468951, 1.045088896507099, 1.05855030690726, 1.064893677016921, 1.041808734830962, 1.051859845029565, 1.048695621970732, 1.062107708436281, 1.068428312874284, 1.020997125708751)
    24 = @(1.02, 1.046901591348246, 1.046327486515711, 1.060256163693427, 1.066807589212314, 1.042270087634172, 1.053018715308494, 1.049661899774238, 1.063543706020096, 1.070073549714957, 1.021394574400787)
    25 = @(1.02, 1.048844982918109, 1.047764097852791, 1.062236967451152, 1.069031103016446, 1.04279966295997, 1.054360374638098, 1.050779496055907, 1.065208666706813, 1.071982732459203, 1.021853998805246)
}

foreach ($r in $data.Keys) {
    $rowvals = $data[$r]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $addr = "{0}{1}" -f $cols[$i], $r
        $ws.Range($addr).Value = $rowvals[$i]
    }
}

Write-Host "applied $($data.Count) rows"
